# Commit: "unify the conception of DataNode, DataTable, Entity."
#
# Content-level changes made in this revision:
#   1. The worksheet formerly called "Property1" is renamed to "DataNode"
#      (part of unifying naming across DataNode/DataTable/Entity sheets).
#   2. The user's active cell/selection on that sheet moved from K17 to B41.
#
# (The rest of the underlying OOXML diff - fileVersion/rupBuild bump,
#  xr/xr2/xr16 revision-tracking namespaces & uids, absPath, window
#  geometry pixels, new phoneticPr/font entry, "Normal"->"常规" cell-style
#  locale text, timeline style extension, default row height / column
#  width micro-drift - is metadata that Excel regenerates by itself when a
#  workbook is opened and resaved by a different Excel build/locale; it is
#  not the product of any user action in the object model, so it is not
#  reproduced here.)

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Property1")

$ws.Name = "DataNode"

$ws.Range("B41").Select()
